$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 6.158193778588456
$ws.Range("C3").Value = 20.16510735538278
$ws.Range("C4").Value = -0.7604140070405252
$ws.Range("C6").Value = 8.31420991993544
$ws.Range("C7").Value = 26.797196171581454
$ws.Range("C8").Value = 0.44695802546495533
$ws.Range("C10").Value = 8.31420991993544
$ws.Range("C11").Value = 26.797196171581454
$ws.Range("C12").Value = 0.44695802546495533
$ws.Range("C14").Value = 7.2761566832357
$ws.Range("C15").Value = 23.604056227168023
$ws.Range("C16").Value = 0.3218170343653928
$ws.Range("C18").Value = 7.205145092317862
$ws.Range("C19").Value = 23.385618540679992
$ws.Range("C20").Value = -0.017311574674248202

# --- WING sheet ---
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C6").Value = 21.255548943013924

# --- FUEL TANK sheet ---
$ws = $wb.Worksheets.Item("FUEL TANK")
$ws.Range("C6").Value = 21.442375947690707

# --- LANDING GEARS sheet ---
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 18.012890076310185
